$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aggregated")

# Rename the header of the first column from "Field Name" to "Column Name"
$ws.Range("B2").Value = "Column Name"

# Remove the stray "Relative Difference" / "Difference" values that were
# accidentally populated for the FIRST_NAME rows of each group (M and Z)
$ws.Range("F3:G3").Clear()
$ws.Range("F10:G10").Clear()
